$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to text so numeric-looking strings like "9.00"
# or "57.20" are not silently coerced to numbers / lose formatting.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '66.838.16'
$ws.Range('E2').Value = '  +2.62%  '
$ws.Range('D3').Value = '3.438.29'
$ws.Range('E3').Value = '  +1.73%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '581.89'
$ws.Range('E5').Value = '  +4.77%  '
$ws.Range('D6').Value = '188.59'
$ws.Range('E6').Value = '  +8.25%  '
$ws.Range('D7').Value = '0.628'
$ws.Range('E7').Value = '  -0.57%  '
$ws.Range('D8').Value = '3.428.86'
$ws.Range('E8').Value = '  +1.80%  '
$ws.Range('E9').Value = '  -0.03%  '
$ws.Range('E10').Value = '  -1.65%  '
$ws.Range('D11').Value = '0.644'
$ws.Range('E11').Value = '  +1.12%  '
$ws.Range('D12').Value = '57.20'
$ws.Range('E12').Value = '  +6.63%  '
$ws.Range('D13').Value = '0.0000276'
$ws.Range('E13').Value = '  -0.89%  '
$ws.Range('D14').Value = '9.44'
$ws.Range('E14').Value = '  +3.12%  '
$ws.Range('D15').Value = '3.990.51'
$ws.Range('E15').Value = '  +1.86%  '
$ws.Range('D16').Value = '18.81'
$ws.Range('E16').Value = '  +2.65%  '
$ws.Range('D17').Value = '3.459.51'
$ws.Range('E17').Value = '  +1.79%  '
$ws.Range('D18').Value = '66.795.09'
$ws.Range('E18').Value = '  +2.75%  '
$ws.Range('B19').Value = 'TRON'
$ws.Range('C19').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D19').Value = '0.119'
$ws.Range('E19').Value = '  +0.33%  '
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').Value = '12.07'
$ws.Range('E20').Value = '  +2.11%  '
$ws.Range('D21').Value = '1.02'
$ws.Range('E21').Value = '  +2.77%  '
$ws.Range('D22').Value = '476.52'
$ws.Range('E22').Value = '  +4.51%  '
$ws.Range('D23').Value = '5.39'
$ws.Range('E23').Value = '  +10.81%  '
$ws.Range('D24').Value = '17.08'
$ws.Range('E24').Value = '  +20.68%  '
$ws.Range('D25').Value = '4.35'
$ws.Range('E25').Value = '  +6.93%  '
$ws.Range('D26').Value = '89.15'
$ws.Range('E26').Value = '  +1.78%  '
$ws.Range('D27').Value = '2.98'
$ws.Range('E27').Value = '  +3.48%  '
$ws.Range('D28').Value = '10.91'
$ws.Range('E28').Value = '  +2.13%  '
$ws.Range('D29').Value = '9.00'
$ws.Range('E29').Value = '  +3.39%  '
$ws.Range('D30').Value = '31.11'
$ws.Range('E30').Value = '  -0.02%  '
$ws.Range('D31').Value = '7.43'
$ws.Range('E31').Value = '  +13.70%  '
$ws.Range('D32').Value = '600.79'
$ws.Range('E32').Value = '  +3.90%  '
$ws.Range('B33').Value = 'OKB'
$ws.Range('C33').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D33').Value = '64.55'
$ws.Range('E33').Value = '  +2.01%  '
$ws.Range('B34').Value = 'Cosmos'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D34').Value = '11.74'
$ws.Range('E34').Value = '  +2.53%  '
$ws.Range('E35').Value = '  +4.02%  '
$ws.Range('E36').Value = '  -0.08%  '
$ws.Range('E37').Value = '  +3.00%  '
$ws.Range('D38').Value = '37.01'
$ws.Range('E38').Value = '  +3.75%  '
$ws.Range('D39').Value = '0.389'
$ws.Range('E39').Value = '  +4.31%  '
$ws.Range('D40').Value = '3.47'
$ws.Range('E40').Value = '  -3.46%  '
$ws.Range('D41').Value = '0.0₃0750'
$ws.Range('E41').Value = '  +1.64%  '
$ws.Range('D42').Value = '3.198.18'
$ws.Range('E42').Value = '  +2.84%  '
$ws.Range('D43').Value = '2.91'
$ws.Range('E43').Value = '  +5.40%  '
$ws.Range('D44').Value = '0.0430'
$ws.Range('E44').Value = '  +3.29%  '
$ws.Range('D45').Value = '2.58'
$ws.Range('E45').Value = '  +5.66%  '
$ws.Range('E46').Value = '  +1.89%  '
$ws.Range('D47').Value = '2.73'
$ws.Range('E47').Value = '  +21.46%  '
$ws.Range('E48').Value = '  +0.62%  '
$ws.Range('D49').Value = '1.00'
$ws.Range('D50').Value = '8.64'
$ws.Range('E50').Value = '  +4.06%  '
$ws.Range('D51').Value = '3.19'
$ws.Range('E51').Value = '  +5.23%  '

# Restore the original (unstyled) look for column D now that the
# text values are locked in - matches the source workbook which has
# no explicit style on these data cells.
$ws.Range("D2:D51").Style = "Normal"

